$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "tfgjui"
$ws.Range("B2").Value = "ansokolkin@gmail.com"
$ws.Range("C2").Value = "tuyijol"

$ws.Range("A3").Value = "ghjkuilok;"
$ws.Range("B3").Value = "ansokolkin@gmail.com"
$ws.Range("C3").Value = "tyuiljo"

$ws.Range("A4").Value = "beta"
$ws.Range("B4").Value = "ansokolkin@gmail.com"
$ws.Range("C4").Value = 'X4#%G$B#y2+RlZYu&Yun1PE'
